# Scheduled runner update: refresh cached market-price figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J32").Value = 634
$ws.Range("H32").Value = 725.5
$ws.Range("L32").Value = 634
$ws.Range("N32").Value = -1286

$ws.Range("J100").Value = 3470.7
$ws.Range("N100").Value = -4552.7
$ws.Range("I100").Value = 1814.6666
$ws.Range("L100").Value = 3470.7
$ws.Range("H100").Value = 2477.08
$ws.Range("M100").Value = -1273.6666
$ws.Range("K100").Value = 1814.6666

$ws.Range("K113").Value = 3936.875
$ws.Range("M113").Value = -682.875
$ws.Range("I113").Value = 3936.875
$ws.Range("H113").Value = 3753.913
$ws.Range("L113").Value = 3335.7144
$ws.Range("J113").Value = 3335.7144
$ws.Range("N113").Value = -9843.714400000001

$ws.Range("K116").Value = 3608.375
$ws.Range("I116").Value = 3608.375
$ws.Range("H116").Value = 4100.7104
$ws.Range("J116").Value = 4944.7144
$ws.Range("N116").Value = -11828.7144
$ws.Range("L116").Value = 4944.7144
$ws.Range("M116").Value = -166.375

$ws.Range("J125").Value = 1295.8572
$ws.Range("L125").Value = 11662.7148
$ws.Range("N125").Value = -16582.7148
$ws.Range("M125").Value = -14028
$ws.Range("I125").Value = 1832
$ws.Range("K125").Value = 16488
$ws.Range("H125").Value = 1597.4375

$ws.Range("N134").Value = -50140
$ws.Range("J134").Value = 40000
$ws.Range("H134").Value = 40000
$ws.Range("L134").Value = 40000

$ws.Range("M141").Value = -5390.734
$ws.Range("K141").Value = 10570.734
$ws.Range("J141").Value = 2086055.2
$ws.Range("I141").Value = 3523.578
$ws.Range("H141").Value = 211776.73
$ws.Range("L141").Value = 6258165.6
$ws.Range("N141").Value = -6268525.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 655
$ws.Range("M32").Value = -814.5834
$ws.Range("H32").Value = 1083.72
$ws.Range("K32").Value = 1101.5834
$ws.Range("I32").Value = 1101.5834
$ws.Range("L32").Value = 655
$ws.Range("N32").Value = -1229

$ws.Range("N132").Value = -13818.2
$ws.Range("L132").Value = 8758.200000000001
$ws.Range("H132").Value = 2050.7273
$ws.Range("K132").Value = 5174.924999999999
$ws.Range("I132").Value = 1724.975
$ws.Range("J132").Value = 2919.4
$ws.Range("M132").Value = -2644.924999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 439.92307
$ws.Range("M64").Value = -93.66665999999998
$ws.Range("K64").Value = 318.66666
$ws.Range("I64").Value = 318.66666

$ws.Range("I67").Value = 318.66666
$ws.Range("M67").Value = 461.33334
$ws.Range("K67").Value = 318.66666
$ws.Range("H67").Value = 439.92307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 3725
$ws.Range("H16").Value = 2862.5
$ws.Range("K16").Value = 2000
$ws.Range("N16").Value = -4299
$ws.Range("M16").Value = -1713
$ws.Range("L16").Value = 3725
$ws.Range("I16").Value = 2000

$ws.Range("J99").Value = 8562.666999999999
$ws.Range("L99").Value = 8562.666999999999
$ws.Range("M99").Value = -3174
$ws.Range("I99").Value = 4672
$ws.Range("K99").Value = 4672
$ws.Range("N99").Value = -11558.667
$ws.Range("H99").Value = 6339.4287

$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170
$ws.Range("I113").Value = 2000
$ws.Range("H113").Value = 2862.5
$ws.Range("L113").Value = 3725
$ws.Range("J113").Value = 3725
$ws.Range("N113").Value = -8065

$ws.Range("N122").Value = -25150
$ws.Range("H122").Value = 4084.4
$ws.Range("L122").Value = 20250
$ws.Range("K122").Value = 9345.2724
$ws.Range("J122").Value = 6750
$ws.Range("M122").Value = -6895.2724
$ws.Range("I122").Value = 3115.0908

$ws.Range("K126").Value = 14016
$ws.Range("N126").Value = -30628.001
$ws.Range("H126").Value = 6339.4287
$ws.Range("L126").Value = 25688.001
$ws.Range("I126").Value = 4672
$ws.Range("J126").Value = 8562.666999999999
$ws.Range("M126").Value = -11546

$ws.Range("H132").Value = 1630.3508
$ws.Range("K132").Value = 3794.487599999999
$ws.Range("I132").Value = 1264.8292
$ws.Range("M132").Value = -1264.487599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K5").Value = 974.4782399999999
$ws.Range("J5").Value = 2238.3333
$ws.Range("N5").Value = -6938.999899999999
$ws.Range("L5").Value = 6714.999899999999
$ws.Range("H5").Value = 720.7241
$ws.Range("I5").Value = 324.82608
$ws.Range("M5").Value = -862.4782399999999

$ws.Range("L135").Value = 20144.9997
$ws.Range("I135").Value = 324.82608
$ws.Range("N135").Value = -25214.9997
$ws.Range("K135").Value = 2923.43472
$ws.Range("H135").Value = 720.7241
$ws.Range("M135").Value = -388.4347199999997
$ws.Range("J135").Value = 2238.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 36700
$ws.Range("L42").Value = 36700
$ws.Range("N42").Value = -37670
$ws.Range("J42").Value = 36700

$ws.Range("N115").Value = -39050
$ws.Range("H115").Value = 36700
$ws.Range("L115").Value = 36700
$ws.Range("J115").Value = 36700

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2599.4285
$ws.Range("M7").Value = -1287.3334
$ws.Range("I7").Value = 1399.3334
$ws.Range("N7").Value = -3723.5
$ws.Range("K7").Value = 1399.3334
$ws.Range("L7").Value = 3499.5
$ws.Range("J7").Value = 3499.5

$ws.Range("I61").Value = 333333980
$ws.Range("H61").Value = 166671330
$ws.Range("K61").Value = 333333980
$ws.Range("M61").Value = -333333778

$ws.Range("K93").Value = 2579.9333
$ws.Range("M93").Value = -1331.9333
$ws.Range("I93").Value = 2579.9333
$ws.Range("H93").Value = 2921.158

$ws.Range("K113").Value = 333333980
$ws.Range("M113").Value = -333331810
$ws.Range("I113").Value = 333333980
$ws.Range("H113").Value = 166671330

$ws.Range("H122").Value = 3266.6667
$ws.Range("K122").Value = 7835.293799999999
$ws.Range("M122").Value = -5385.293799999999
$ws.Range("I122").Value = 2611.7646

$ws.Range("K126").Value = 4198.0002
$ws.Range("N126").Value = -15438.5
$ws.Range("H126").Value = 2599.4285
$ws.Range("L126").Value = 10498.5
$ws.Range("I126").Value = 1399.3334
$ws.Range("J126").Value = 3499.5
$ws.Range("M126").Value = -1728.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J100").Value = 1980
$ws.Range("N100").Value = -5042
$ws.Range("I100").Value = 2000
$ws.Range("L100").Value = 3960
$ws.Range("H100").Value = 1983.6364
$ws.Range("M100").Value = -3459
$ws.Range("K100").Value = 4000

$ws.Range("K126").Value = 3897.84
$ws.Range("N126").Value = -37515971
$ws.Range("H126").Value = 3032178.8
$ws.Range("L126").Value = 37511031
$ws.Range("I126").Value = 1299.28
$ws.Range("J126").Value = 12503677
$ws.Range("M126").Value = -1427.84

$ws.Range("L135").Value = 66611.25
$ws.Range("N135").Value = -76751.25
$ws.Range("H135").Value = 66611.25
$ws.Range("J135").Value = 66611.25
